# datatest_CRM edit: add getExcelData-style extra row to the Login sheet,
# clean up the now-redundant "user@gmail.com" sample, and re-point the
# active tab / selections to match the refreshed workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Login
$ws2 = $wb.Worksheets.Item(2)   # Customer

function Remove-HyperlinksAt($ws, [string[]]$targets) {
    $changed = $true
    while ($changed) {
        $changed = $false
        foreach ($h in $ws.Hyperlinks) {
            $addr = $h.Range.Address()
            if ($targets -contains $addr) {
                $h.Delete()
                $changed = $true
                break
            }
        }
    }
}

# --- Login sheet --------------------------------------------------------

# Row 3 used to hold the throw-away "user@gmail.com" sample; it now mirrors
# row 2's admin account instead (the unused shared string drops out on save).
$ws1.Range("A3").Value = "admin@example.com"

# Rows 2 & 3 no longer carry a live hyperlink -- drop the relationships,
# then strip the inherited Hyperlink font (underline/colour) while keeping
# the cells left-aligned.
Remove-HyperlinksAt $ws1 @('$A$2', '$A$3')

$r2 = $ws1.Range("A2")
$r2.HorizontalAlignment = -4131   # xlLeft
$r2.Font.Underline = 0

$r3 = $ws1.Range("A3")
$r3.HorizontalAlignment = -4131   # xlLeft
$r3.Font.Underline = 0

# New row 4: another login sample, this time a real hyperlink cell again.
$ws1.Range("A4").Value = "admin1@example.com"
$ws1.Range("B4").Value = 123456

$ws1.Hyperlinks.Add($ws1.Range("A4"), "mailto:admin1@example.com")
$ws1.Range("A4").HorizontalAlignment = -4131   # xlLeft

# Print setup now explicitly set to portrait.
$ws1.PageSetup.Orientation = 1   # xlPortrait

# --- Customer sheet selection -------------------------------------------

$ws2.Activate()
$ws2.Range("A8").Select()

# --- Active tab switches back to Login -----------------------------------

$ws1.Activate()
$ws1.Range("B6").Select()
